# This workbook's single sheet contains a table of price records for
# "Zanahoria" (carrots) at the "Vega Modelo de Temuco" market. A new daily
# record needs to be inserted as a new row at position 185, pushing all
# the existing rows from 185 downward down by one (so the former row 266
# becomes row 267).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185; this shifts rows 185:266 down to 186:267 and
# grows the sheet dimension from A1:R266 to A1:R267 automatically.
$ws.Rows("185:185").Insert()

# Populate the newly inserted row 185 with the new record. The
# "constant" identifying columns (Mercado ID, Mercado, Region, Codreg,
# Categoria ID, Categoria, Variedad, Calidad, Clasificacion) match every
# other row in this table.
$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 44609
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = 100114013
$ws.Range("G185").Value = "Zanahoria"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 80
$ws.Range("K185").Value = 8000
$ws.Range("L185").Value = 8000
$ws.Range("M185").Value = 8000
$ws.Range("N185").Value = "$/saco 25 kilos"
$ws.Range("O185").Value = "Provincia de Cautín"
$ws.Range("P185").Value = 320
$ws.Range("Q185").Value = 25
$ws.Range("R185").Value = "Hortaliza"
